$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 3.4
$ws.Range("S3").Value = 1.3
$ws.Range("T3").Value = 3.4
$ws.Range("W3").Value = 11
$ws.Range("AT3").Value = 3.4
$ws.Range("AY3").Value = 21
